$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.560.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.935.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.94"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4833"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2908"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06808"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "111.95"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.47"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.933.74"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07580"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.469"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6798"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "298.85"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.567.78"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007653"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.585"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.187.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9974"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.495"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.510"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.45"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.46"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.148"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1070"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.437"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.158"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.093"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04989"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7428"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.149"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02036"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.691"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.033"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "109.95"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4465"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8703"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.835"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.95%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.75"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.84%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.283"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.22"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.346"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1236"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.2535"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.99"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.56%  "
